# Refactored parse tree builder to take into account optimised FUNCTIONS
# that use skip tokens - add new test cells exercising IF() with skipped
# (optional) arguments.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("ExcelFormulaTreeTest")
$wsLookup = $wb.Worksheets.Item("Lookup")

# New row 12: boolean IF() with all three arguments supplied.
$wsMain.Range("E12").Value = 50
$wsMain.Range("C12").Value = $true
$wsMain.Range("D12").Value = $false
$wsMain.Range("A12").Formula = "=IF(E12>(1+50),C12,D12)"

# New row 13: IF() with the optional third argument skipped.
$wsMain.Range("A13").Formula = "=IF(E12>0,""Yes"")"

# Update selections to match the new state left behind on each sheet.
$wsLookup.Activate()
$wsLookup.Range("F17").Select()

$wsMain.Activate()
$wsMain.Range("C12:E12").Select()
